$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 12990691
$ws.Range("I74").Value = 12990691
$ws.Range("K74").Value = 12990691
$ws.Range("M74").Value = -12989755
$ws.Range("H77").Value = 12990691
$ws.Range("I77").Value = 12990691
$ws.Range("K77").Value = 64953455
$ws.Range("M77").Value = -64948775
$ws.Range("H97").Value = 2646.7
$ws.Range("J97").Value = 2829.6667
$ws.Range("L97").Value = 8489.000100000001
$ws.Range("N97").Value = -9481.000100000001
$ws.Range("H137").Value = 2309.5417
$ws.Range("I137").Value = 1244.05
$ws.Range("K137").Value = 3732.15
$ws.Range("M137").Value = -1182.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2731.117
$ws.Range("I2").Value = 1272.4193
$ws.Range("J2").Value = 8760.4
$ws.Range("K2").Value = 1272.4193
$ws.Range("L2").Value = 8760.4
$ws.Range("M2").Value = -1159.4193
$ws.Range("N2").Value = -8986.4
$ws.Range("H32").Value = 1548398.1
$ws.Range("I32").Value = 4242.646
$ws.Range("K32").Value = 4242.646
$ws.Range("M32").Value = -3955.646
$ws.Range("H61").Value = 8763.806
$ws.Range("I61").Value = 4893.9688
$ws.Range("K61").Value = 4893.9688
$ws.Range("M61").Value = -4681.9688
$ws.Range("H74").Value = 6260.9
$ws.Range("I74").Value = 7011.636
$ws.Range("K74").Value = 7011.636
$ws.Range("M74").Value = -6137.636
$ws.Range("H77").Value = 6260.9
$ws.Range("I77").Value = 7011.636
$ws.Range("K77").Value = 35058.18
$ws.Range("M77").Value = -30690.18
$ws.Range("H116").Value = 2731.117
$ws.Range("I116").Value = 1272.4193
$ws.Range("J116").Value = 8760.4
$ws.Range("K116").Value = 1272.4193
$ws.Range("L116").Value = 8760.4
$ws.Range("M116").Value = 1021.5807
$ws.Range("N116").Value = -13348.4
$ws.Range("H132").Value = 660620.7
$ws.Range("I132").Value = 734021.25
$ws.Range("K132").Value = 2202063.75
$ws.Range("M132").Value = -2199533.75
$ws.Range("H136").Value = 8763.806
$ws.Range("I136").Value = 4893.9688
$ws.Range("K136").Value = 14681.9064
$ws.Range("M136").Value = -12131.9064

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2731.117
$ws.Range("I3").Value = 1272.4193
$ws.Range("J3").Value = 8760.4
$ws.Range("K3").Value = 1272.4193
$ws.Range("L3").Value = 8760.4
$ws.Range("M3").Value = -1158.4193
$ws.Range("N3").Value = -8988.4
$ws.Range("H86").Value = 4647.174
$ws.Range("I86").Value = 2949.4167
$ws.Range("J86").Value = 6499.273
$ws.Range("K86").Value = 2949.4167
$ws.Range("L86").Value = 6499.273
$ws.Range("M86").Value = -1826.4167
$ws.Range("N86").Value = -8745.273000000001
$ws.Range("H89").Value = 4647.174
$ws.Range("I89").Value = 2949.4167
$ws.Range("J89").Value = 6499.273
$ws.Range("K89").Value = 14747.0835
$ws.Range("L89").Value = 32496.365
$ws.Range("M89").Value = -9131.083500000001
$ws.Range("N89").Value = -43728.36500000001
$ws.Range("H99").Value = 7701.909
$ws.Range("I99").Value = 7554.5557
$ws.Range("J99").Value = 7909.125
$ws.Range("K99").Value = 7554.5557
$ws.Range("L99").Value = 7909.125
$ws.Range("M99").Value = -6056.5557
$ws.Range("N99").Value = -10905.125
$ws.Range("H134").Value = 758814.9
$ws.Range("I134").Value = 869257.9399999999
$ws.Range("J134").Value = 7802.2
$ws.Range("K134").Value = 2607773.82
$ws.Range("L134").Value = 23406.6
$ws.Range("M134").Value = -2605238.82
$ws.Range("N134").Value = -28476.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 620.8889
$ws.Range("I5").Value = 311.42856
$ws.Range("J5").Value = 1704
$ws.Range("K5").Value = 311.42856
$ws.Range("L5").Value = 1704
$ws.Range("M5").Value = -199.42856
$ws.Range("N5").Value = -1928
$ws.Range("H31").Value = 7964.15
$ws.Range("I31").Value = 9165
$ws.Range("K31").Value = 9165
$ws.Range("M31").Value = -8870
$ws.Range("H34").Value = 7964.15
$ws.Range("I34").Value = 9165
$ws.Range("K34").Value = 9165
$ws.Range("M34").Value = -8963
$ws.Range("H47").Value = 15977.75
$ws.Range("I47").Value = 10000
$ws.Range("J47").Value = 17970.334
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 17970.334
$ws.Range("M47").Value = -9434
$ws.Range("N47").Value = -19102.334
$ws.Range("H59").Value = 34130.223
$ws.Range("J59").Value = 34646.5
$ws.Range("L59").Value = 34646.5
$ws.Range("N59").Value = -36936.5
$ws.Range("H60").Value = 24646.357
$ws.Range("J60").Value = 45826
$ws.Range("L60").Value = 45826
$ws.Range("N60").Value = -46848
$ws.Range("H141").Value = 226939
$ws.Range("J141").Value = 245132.9
$ws.Range("L141").Value = 245132.9
$ws.Range("N141").Value = -255492.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 29000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 90000
$ws.Range("N87").Value = -92496
$ws.Range("H90").Value = 29000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 270000
$ws.Range("N90").Value = -282480
$ws.Range("H103").Value = 949
$ws.Range("I103").Value = 786
$ws.Range("K103").Value = 2358
$ws.Range("M103").Value = -1479
$ws.Range("H107").Value = 3029.1592
$ws.Range("I107").Value = 550.25
$ws.Range("J107").Value = 3277.05
$ws.Range("K107").Value = 1650.75
$ws.Range("L107").Value = 9831.150000000001
$ws.Range("M107").Value = 269.25
$ws.Range("N107").Value = -13671.15
$ws.Range("H117").Value = 3281.111
$ws.Range("I117").Value = 855
$ws.Range("J117").Value = 8133.3335
$ws.Range("K117").Value = 2565
$ws.Range("L117").Value = 24400.0005
$ws.Range("M117").Value = 877
$ws.Range("N117").Value = -31284.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5831.5454
$ws.Range("I80").Value = 3199.7144
$ws.Range("J80").Value = 10437.25
$ws.Range("K80").Value = 3199.7144
$ws.Range("L80").Value = 10437.25
$ws.Range("M80").Value = -2201.7144
$ws.Range("N80").Value = -12433.25
$ws.Range("H83").Value = 5831.5454
$ws.Range("I83").Value = 3199.7144
$ws.Range("J83").Value = 10437.25
$ws.Range("K83").Value = 15998.572
$ws.Range("L83").Value = 52186.25
$ws.Range("M83").Value = -11006.572
$ws.Range("N83").Value = -62170.25
$ws.Range("H122").Value = 4037.0312
$ws.Range("I122").Value = 3136.9546
$ws.Range("J122").Value = 6017.2
$ws.Range("K122").Value = 9410.863799999999
$ws.Range("L122").Value = 18051.6
$ws.Range("M122").Value = -6960.863799999999
$ws.Range("N122").Value = -22951.6
$ws.Range("H126").Value = 41682284
$ws.Range("I126").Value = 166669340
$ws.Range("J126").Value = 19932
$ws.Range("K126").Value = 500008020
$ws.Range("L126").Value = 59796
$ws.Range("M126").Value = -500005550
$ws.Range("N126").Value = -64736
$ws.Range("H132").Value = 5547.265
$ws.Range("I132").Value = 5281.533
$ws.Range("J132").Value = 5966.8423
$ws.Range("K132").Value = 15844.599
$ws.Range("L132").Value = 17900.5269
$ws.Range("M132").Value = -13314.599
$ws.Range("N132").Value = -22960.5269

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 26317688
$ws.Range("I46").Value = 1069.6666
$ws.Range("J46").Value = 38463820
$ws.Range("K46").Value = 1069.6666
$ws.Range("L46").Value = 38463820
$ws.Range("M46").Value = -881.6666
$ws.Range("N46").Value = -38464196
$ws.Range("H61").Value = 7491.6294
$ws.Range("I61").Value = 5337.087
$ws.Range("J61").Value = 19880.25
$ws.Range("K61").Value = 5337.087
$ws.Range("L61").Value = 19880.25
$ws.Range("M61").Value = -5135.087
$ws.Range("N61").Value = -20284.25
$ws.Range("H68").Value = 2195.5
$ws.Range("J68").Value = 1823
$ws.Range("L68").Value = 1823
$ws.Range("N68").Value = -3321
$ws.Range("H71").Value = 2195.5
$ws.Range("J71").Value = 1823
$ws.Range("L71").Value = 9115
$ws.Range("N71").Value = -16603
$ws.Range("H113").Value = 7491.6294
$ws.Range("I113").Value = 5337.087
$ws.Range("J113").Value = 19880.25
$ws.Range("K113").Value = 5337.087
$ws.Range("L113").Value = 19880.25
$ws.Range("M113").Value = -3167.087
$ws.Range("N113").Value = -24220.25
$ws.Range("H136").Value = 15901.8
$ws.Range("I136").Value = 20799.666
$ws.Range("J136").Value = 13802.714
$ws.Range("K136").Value = 62398.99800000001
$ws.Range("L136").Value = 41408.142
$ws.Range("M136").Value = -59848.99800000001
$ws.Range("N136").Value = -46508.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 19995
$ws.Range("J43").Value = 19995
$ws.Range("L43").Value = 19995
$ws.Range("N43").Value = -20293
$ws.Range("H109").Value = 41555
$ws.Range("J109").Value = 41555
$ws.Range("L109").Value = 41555
$ws.Range("N109").Value = -44329
$ws.Range("H132").Value = 5331.8984
$ws.Range("I132").Value = 4754.125
$ws.Range("J132").Value = 6548.263
$ws.Range("K132").Value = 14262.375
$ws.Range("L132").Value = 19644.789
$ws.Range("M132").Value = -11732.375
$ws.Range("N132").Value = -24704.789
$ws.Range("H136").Value = 10426172
$ws.Range("I136").Value = 29428382
$ws.Range("J136").Value = 5605.2905
$ws.Range("K136").Value = 88285146
$ws.Range("L136").Value = 16815.8715
$ws.Range("M136").Value = -88282596
$ws.Range("N136").Value = -21915.8715
